# The two sighting records that were entered on sheet rows 16 and 17 had
# their data swapped: what used to be reported on row 16 moved to row 17,
# and vice versa. The physical rows (and their handful of shared
# formatting) stay where they are - only the field values move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold a value on both rows 16 and 17 and are safe to swap
# with a plain value assignment (no risk of Excel reinterpreting the
# string as a date/number, since these are plain numerics or plain
# words).
$simpleCols = @("A","B","E","F","G","H","P","Q","R","S","AW","AX")

foreach ($col in $simpleCols) {
    $r16 = $ws.Range($col + "16")
    $r17 = $ws.Range($col + "17")

    $v16 = $r16.Value2
    $v17 = $r17.Value2

    $r16.Value2 = $v17
    $r17.Value2 = $v16
}

# Columns K, L, N only carried an (empty) cell on row 17; row 16 had no
# cell there at all. After the swap row 16 should gain those empty
# cells and row 17 should lose them. Range.Copy preserves "present but
# empty" cells (unlike assigning Value2 = ""), so copy row 17's empty
# cells onto row 16, then clear them from row 17.
$ws.Range("K17:L17").Copy($ws.Range("K16"))
$ws.Range("N17").Copy($ws.Range("N16"))
$ws.Range("K17:L17").ClearContents()
$ws.Range("N17").ClearContents()

# Column M held "äldre spår" on row 17 and nothing on row 16; column AC
# held "Ringhack på gran" on row 17 and nothing on row 16. Move the text
# to row 16 and remove it from row 17 entirely.
$ws.Range("M16").Value2 = $ws.Range("M17").Value2
$ws.Range("M17").ClearContents()

$ws.Range("AC16").Value2 = $ws.Range("AC17").Value2
$ws.Range("AC17").ClearContents()
